# #327 Ajout des profils d'acces
# - Regenerated "Date" metadata value (Metadata!B8)
# - Swapped the two mapping columns ("Mapping: RIM Mapping" and
#   "Mapping: Spécification métier vers l'extension ROR ClosingType")
#   on the Elements sheet, including their header cells, their data and
#   their column widths (columns AK/AL, i.e. 37/38).

$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: refresh the generation Date ---------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2024-03-19T13:17:15+00:00"

# --- 2. Elements sheet: swap columns AK (37) and AL (38) ----------------
$wsElem = $wb.Worksheets.Item("Elements")

$usedRows = $wsElem.UsedRange.Rows.Count
for ($r = 1; $r -le $usedRows; $r++) {
    $akCell = $wsElem.Cells.Item($r, 37)
    $alCell = $wsElem.Cells.Item($r, 38)

    $akVal = $akCell.Value2
    $alVal = $alCell.Value2

    $akCell.Value = $alVal
    $alCell.Value = $akVal
}

# Swap the (best-fit) column widths that travelled with the data.
$wsElem.Columns.Item(37).ColumnWidth = 67.61328125
$wsElem.Columns.Item(38).ColumnWidth = 24.98046875
